$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New detection rows to append (Fly_ID, Class, First_Detection_Date (serial),
# First_Detection_Image, Placa ID, Localização, Latitude, Longitude,
# First_Coords, First_Confidence)
$newRows = @(
    @("312b9e21-1bbf-4ba0-ad02-ea0aaa6f842b", "mosca", 45876, "image_20250807111026_ppp0.jpg", "PLACA_20250717165933", "Beja", 38.02035, -7.94715, "641,529,688,576", "0.75"),
    @("d46143c3-a87f-42fc-a10d-2d5a22dc0c2b", "mosca", 45876, "image_20250807111026_ppp0.jpg", "PLACA_20250717165933", "Beja", 38.02035, -7.94715, "793,481,831,526", "0.70")
)

$startRow = 16

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    $dateCell = $ws.Cells.Item($r, 3)
    $dateCell.Value = $row[2]
    $dateCell.Style = $ws.Cells.Item(2, 3).Style
    $dateCell.NumberFormat = $ws.Cells.Item(2, 3).NumberFormat

    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]

    # Columns I (bounding-box coords) and J (confidence) are stored as plain
    # text in the source data (e.g. "641,529,688,576", "0.75") even though
    # they look numeric. Force text entry (like Excel's Text-formatted
    # cells) so commas/decimals aren't reinterpreted as a number, then drop
    # back to the default style once the literal text is committed.
    $coordCell = $ws.Cells.Item($r, 9)
    $coordCell.NumberFormat = "@"
    $coordCell.Value = $row[8]
    $coordCell.Style = "Normal"

    $confCell = $ws.Cells.Item($r, 10)
    $confCell.NumberFormat = "@"
    $confCell.Value = $row[9]
    $confCell.Style = "Normal"
}
